# Update BOM: add 40mm stepper motors and GT2 timing belt pulley rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "40mm stepper motors"
$ws.Range("B8").Value = "x2"
$ws.Hyperlinks.Add($ws.Range("C8"), "https://www.amazon.nl/-/en/gp/product/B07SQNYZDY/ref=ppx_yo_dt_b_search_asin_title?ie=UTF8&th=1", "", "", "https://www.amazon.nl/-/en/gp/product/B07SQNYZDY/ref=ppx_yo_dt_b_search_asin_title?ie=UTF8&th=1")

$ws.Range("A9").Value = "GT2 Timing Belt Pulley 20 tooth, 6mm width, 5mm bore"
$ws.Range("B9").Value = "x2"
$ws.Hyperlinks.Add($ws.Range("C9"), "https://www.aliexpress.com/item/10000046456013.html?spm=a2g0o.productlist.main.5.59ba2648jXfmXx&algo_pvid=3f3436a0-e6bd-49ab-8fd0-6cb43f8f0c6e&aem_p4p_detail=2023041223572512566913022360240000127032&algo_exp_id=3f3436a0-e6bd-49ab-8fd0-6cb43f8f0c6e-2&pdp_npi=3%40dis%21EUR%212.06%211.44%21%21%21%21%21%402145279016813690456322177d0715%2112000025024121597%21sea%21NL%212329305485&curPageLogUid=09OYRIFzKRYK&ad_pvid=2023041223572512566913022360240000127032_3&ad_pvid=2023041223572512566913022360240000127032_3", "", "", "https://www.aliexpress.com/item/10000046456013.html?spm=a2g0o.productlist.main.5.59ba2648jXfmXx&algo_pvid=3f3436a0-e6bd-49ab-8fd0-6cb43f8f0c6e&aem_p4p_detail=2023041223572512566913022360240000127032&algo_exp_id=3f3436a0-e6bd-49ab-8fd0-6cb43f8f0c6e-2&pdp_npi=3%40dis%21EUR%212.06%211.44%21%21%21%21%21%402145279016813690456322177d0715%2112000025024121597%21sea%21NL%212329305485&curPageLogUid=09OYRIFzKRYK&ad_pvid=2023041223572512566913022360240000127032_3&ad_pvid=2023041223572512566913022360240000127032_3")

$ws.Columns.Item(1).ColumnWidth = 49.24

$ws.Range("A9").Select()
